$p = $ppt.ActivePresentation

# The "auto slide number alignment" fix materializes explicit geometry /
# body / list-style formatting (copied down from the slide master's own
# "Slide Number Placeholder" definition) onto the inline slide-number
# placeholder shapes that sit on slide positions 2 and 6 (the two slides
# that were touched by the slide-2-to-6 swap).
$targetPositions = @(2, 6)

foreach ($pos in $targetPositions) {
    $s = $p.Slides.Item($pos)

    # The slide number placeholder is always the first shape on these slides.
    $sh = $s.Shapes.Item(1)

    # Rename to match the fixed-up placeholder name.
    $sh.Name = "Slide Number Placeholder 1"

    # Materialize the inherited position/size as an explicit <a:xfrm>.
    $sh.Left = 516
    $sh.Top = 500.5
    $sh.Width = 168
    $sh.Height = 28.75

    # Materialize an explicit rectangle preset geometry <a:prstGeom prst="rect">.
    $sh.AutoShapeType = 1

    $tf = $sh.TextFrame

    # Materialize the inherited text-box body properties.
    $tf.Orientation = 1
    $tf.MarginLeft = 7.2
    $tf.MarginTop = 3.6
    $tf.MarginRight = 7.2
    $tf.MarginBottom = 3.6
    $tf.VerticalAnchor = 3
}
